# Apply the target edit: remove the "Amen Thompson" row and reorder the
# remaining player rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data (Oyuncu Adı, Pozisyon, Takım) in the desired order.
$data = @(
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets")
)

# Clear the old data area (below header) including the row that will be
# removed, then write the new ordering back in.
$lastRow = 18
$ws.Range("A2:C$lastRow").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
